# "Generate Report for Handback"
# Refresh the handback report timestamps / status that get regenerated
# whenever the report is produced again.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-10-17 17:13:45"
$wsOverview.Range("G3").Value = "2016-10-17 17:13:45"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Priority column (E): ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
# Correspond Handoff Datetime column (H)
$wsZhCn.Range("H2").Value = "2016-10-17 17:13:22"
$wsZhCn.Range("H3").Value = "2016-10-17 17:13:22"
# Correspond Handback DateTime column (K)
$wsZhCn.Range("K2").Value = "2016-10-17 17:14:30"
$wsZhCn.Range("K3").Value = "2016-10-17 17:14:30"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# Priority column (E): ht -> mt
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
# Correspond Handback DateTime column (K)
$wsDeDe.Range("K2").Value = "2016-10-17 17:15:14"
$wsDeDe.Range("K3").Value = "2016-10-17 17:15:14"
